$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on column E (the 5th column of the A1:E29 table) that
# shows only rows where column E is blank. This hides every row whose
# "Analise" column contains "-" (rows 3, 10, 11, 13, 14, 15, 16, 21, 24, 25),
# matching Excel's "(Blanks)" filter checkbox behavior.
$rng = $ws.Range("A1:E29")
$rng.AutoFilter(5, @(""))

# Reflect the resulting multi-row selection left over from filtering.
$ws.Range("A2:A29").Select()
